# Automatische test-sync: 2025-08-13 20:15:50
# Append the new mail-log entry to the "Logs" sheet and refresh the
# "Dashboard" summary count to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# New row 5 values (as plain text, matching the existing inline-string rows).
$logs.Cells.Item(5, 1).Value = "Demo inplannen"
$logs.Cells.Item(5, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(5, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(5, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(5, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(5, 6).Value = "2025-08-13 20:15:02"
$logs.Cells.Item(5, 7).Value = "Nee"
$logs.Cells.Item(5, 8).Value = "Ja"
$logs.Cells.Item(5, 9).Value = "Nee"
$logs.Cells.Item(5, 10).Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too.
$ranges = @("D2:D4", "G2:G4", "H2:H4", "I2:I4", "J2:J4")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newRange = $logs.Range("$($col)2:$($col)5")
    $fcs = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for this category.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 4
